$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.862.34'
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.410.92'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.34'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.06'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '8.03'
$ws.Range('E9').Value = '  +4.51%  '
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.412'
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.996.65'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.14'
$ws.Range('E14').Value = '  -5.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.427.27'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.927.53'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.35'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.44'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.86'
$ws.Range('E20').Value = '  -4.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.54'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.565'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '75.12'
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.558.17'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000110'
$ws.Range('E26').Value = '  -4.88%  '
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.59'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.89'
$ws.Range('E30').Value = '  -3.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.11'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.00'
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.45'
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '169.67'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.86'
$ws.Range('E38').Value = '  -3.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '30.99'
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.446.68'
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0770'
$ws.Range('E41').Value = '  +0.90%  '
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.35'
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.15'
$ws.Range('E46').Value = '  -5.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.535.59'
$ws.Range('E47').Value = '  -3.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.45'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  -4.99%  '
